$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(716).Insert()

$ws.Cells.Item(716, 1).Value = 3
$ws.Cells.Item(716, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(716, 3).Value = "Coquimbo"
$ws.Cells.Item(716, 4).Value = 45212
$ws.Cells.Item(716, 5).Value = 5
$ws.Cells.Item(716, 6).Value = 100112021
$ws.Cells.Item(716, 7).Value = "Ají"
$ws.Cells.Item(716, 8).Value = "Inferno"
$ws.Cells.Item(716, 9).Value = "Primera"
$ws.Cells.Item(716, 10).Value = 48
$ws.Cells.Item(716, 11).Value = 30000
$ws.Cells.Item(716, 12).Value = 31000
$ws.Cells.Item(716, 13).Value = 30208
$ws.Cells.Item(716, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(716, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(716, 16).Value = 3021
$ws.Cells.Item(716, 17).Value = 10
$ws.Cells.Item(716, 18).Value = "Hortaliza"
